$wb = $excel.ActiveWorkbook

# --- Content change -------------------------------------------------------
# "Ready for handoff" -> "In Translation" everywhere it appears:
#   Overview!E2 (zh-cn status), Overview!F2 (de-de status),
#   zh-cn!C2 (Status), de-de!C2 (Status)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width change ---------------------------------------------------
# The regenerated report narrowed the status columns (their content got
# shorter: "Ready for handoff" -> "In Translation"), so the columns that
# display/hold that text get resized to the new, narrower fitted width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5013020833333
$wsOverview.Columns.Item(6).ColumnWidth = 12.5013020833333
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5013020833333
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5013020833333
